$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric "Qty executed upto date" (column C) updates ---
$ws.Range("C8").Value = 67
$ws.Range("C9").Value = 8
$ws.Range("C10").Value = 65
$ws.Range("C11").Value = 82
$ws.Range("C12").Value = 26
$ws.Range("C13").Value = 73
$ws.Range("C14").Value = 58
$ws.Range("C15").Value = 7
$ws.Range("C16").Value = 67
$ws.Range("C17").Value = 83

# --- Amount columns (G/H) are stored as text ("12345.67") rather than
#     numbers in the workbook, so a plain .Value assignment of a numeric
#     looking string would be auto-coerced to a Number by Excel's normal
#     smart-entry heuristics. Use a literal string formula and convert it
#     to a static value via copy / paste-special so the saved cell keeps
#     its original text type (no residual formula, no style change). ---
function Set-TextValue([string]$addr, [string]$text) {
    $cell = $ws.Range($addr)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue "G9"  "2048.00"
Set-TextValue "G10" "30680.00"
Set-TextValue "G11" "54284.00"
Set-TextValue "G13" "9928.00"
Set-TextValue "G14" "1334.00"

Set-TextValue "G19" "98274.00"
Set-TextValue "H19" "98274.00"
Set-TextValue "G21" "98274.00"
Set-TextValue "H21" "98274.00"

$excel.CutCopyMode = $false
